$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 601
$ws.Range("I29").Value = 601
$ws.Range("K29").Value = 1803
$ws.Range("M29").Value = -1522
# Row 38
$ws.Range("H38").Value = 2025.8113
$ws.Range("J38").Value = 3904.08
$ws.Range("L38").Value = 11712.24
$ws.Range("N38").Value = -12456.24
# Row 74
$ws.Range("H74").Value = 4333.3335
$ws.Range("I74").Value = 4250
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 4250
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -3314
$ws.Range("N74").Value = -6372
# Row 76
$ws.Range("H76").Value = 6176195
$ws.Range("I76").Value = 9262009
$ws.Range("K76").Value = 9262009
$ws.Range("M76").Value = -9261694
# Row 77
$ws.Range("H77").Value = 4333.3335
$ws.Range("I77").Value = 4250
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 21250
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -16570
$ws.Range("N77").Value = -31860
# Row 79
$ws.Range("H79").Value = 6176195
$ws.Range("I79").Value = 9262009
$ws.Range("K79").Value = 9262009
$ws.Range("M79").Value = -9260917
# Row 132
$ws.Range("H132").Value = 1913.8438
$ws.Range("I132").Value = 1780.88
$ws.Range("J132").Value = 2388.7144
$ws.Range("K132").Value = 5342.64
$ws.Range("L132").Value = 7166.1432
$ws.Range("M132").Value = -2812.64
$ws.Range("N132").Value = -12226.1432
# Row 135
$ws.Range("H135").Value = 3566.8386
$ws.Range("I135").Value = 2590.5715
$ws.Range("K135").Value = 23315.1435
$ws.Range("M135").Value = -20780.1435
# Row 137
$ws.Range("H137").Value = 2242.6
$ws.Range("I137").Value = 1419.125
$ws.Range("J137").Value = 2791.5833
$ws.Range("K137").Value = 4257.375
$ws.Range("L137").Value = 8374.749899999999
$ws.Range("M137").Value = -1707.375
$ws.Range("N137").Value = -13474.7499
# Row 138
$ws.Range("H138").Value = 4588.2627
$ws.Range("I138").Value = 2939.9355
$ws.Range("J138").Value = 5339.706
$ws.Range("K138").Value = 8819.806500000001
$ws.Range("L138").Value = 16019.118
$ws.Range("M138").Value = -3679.806500000001
$ws.Range("N138").Value = -26299.118

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1072125.6
$ws.Range("I32").Value = 11972.25
$ws.Range("J32").Value = 43478260
$ws.Range("K32").Value = 11972.25
$ws.Range("L32").Value = 43478260
$ws.Range("M32").Value = -11685.25
$ws.Range("N32").Value = -43478834
# Row 37
$ws.Range("H37").Value = 6682.6
$ws.Range("J37").Value = 8228.25
$ws.Range("L37").Value = 8228.25
$ws.Range("N37").Value = -8774.25
# Row 61
$ws.Range("H61").Value = 4396.6665
$ws.Range("I61").Value = 4422.364
$ws.Range("J61").Value = 4114
$ws.Range("K61").Value = 4422.364
$ws.Range("L61").Value = 4114
$ws.Range("M61").Value = -4210.364
$ws.Range("N61").Value = -4538
# Row 74
$ws.Range("H74").Value = 2249.0908
$ws.Range("I74").Value = 2373.3
$ws.Range("K74").Value = 2373.3
$ws.Range("M74").Value = -1499.3
# Row 77
$ws.Range("H77").Value = 2249.0908
$ws.Range("I77").Value = 2373.3
$ws.Range("K77").Value = 11866.5
$ws.Range("M77").Value = -7498.5
# Row 132
$ws.Range("H132").Value = 21186.963
$ws.Range("I132").Value = 2108.1777
$ws.Range("K132").Value = 6324.533100000001
$ws.Range("M132").Value = -3794.533100000001
# Row 136
$ws.Range("H136").Value = 4396.6665
$ws.Range("I136").Value = 4422.364
$ws.Range("J136").Value = 4114
$ws.Range("K136").Value = 13267.092
$ws.Range("L136").Value = 12342
$ws.Range("M136").Value = -10717.092
$ws.Range("N136").Value = -17442

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 50
$ws.Range("H50").Value = 17632.857
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 17632.857
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 17632.857
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -18780.857

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9946.598
$ws.Range("I31").Value = 3148.275
$ws.Range("J31").Value = 18444.5
$ws.Range("K31").Value = 3148.275
$ws.Range("L31").Value = 18444.5
$ws.Range("M31").Value = -2853.275
$ws.Range("N31").Value = -19034.5
# Row 34
$ws.Range("H34").Value = 9946.598
$ws.Range("I34").Value = 3148.275
$ws.Range("J34").Value = 18444.5
$ws.Range("K34").Value = 3148.275
$ws.Range("L34").Value = 18444.5
$ws.Range("M34").Value = -2946.275
$ws.Range("N34").Value = -18848.5
# Row 51
$ws.Range("H51").Value = 8682
$ws.Range("J51").Value = 9418.4
$ws.Range("L51").Value = 9418.4
$ws.Range("N51").Value = -10890.4
# Row 59
$ws.Range("H59").Value = 12517.667
$ws.Range("J59").Value = 12517.667
$ws.Range("L59").Value = 12517.667
$ws.Range("N59").Value = -14807.667
# Row 60
$ws.Range("H60").Value = 6888.125
$ws.Range("J60").Value = 8278
$ws.Range("L60").Value = 8278
$ws.Range("N60").Value = -9300
# Row 61
$ws.Range("H61").Value = 8682
$ws.Range("J61").Value = 9418.4
$ws.Range("L61").Value = 9418.4
$ws.Range("N61").Value = -10114.4
# Row 68
$ws.Range("H68").Value = 17289
$ws.Range("J68").Value = 17289
$ws.Range("L68").Value = 17289
$ws.Range("N68").Value = -18787
# Row 71
$ws.Range("H71").Value = 17289
$ws.Range("J71").Value = 17289
$ws.Range("L71").Value = 51867
$ws.Range("N71").Value = -59355
# Row 74
$ws.Range("H74").Value = 17995
$ws.Range("J74").Value = 17995
$ws.Range("L74").Value = 17995
$ws.Range("N74").Value = -19743
# Row 77
$ws.Range("H77").Value = 17995
$ws.Range("J77").Value = 17995
$ws.Range("L77").Value = 53985
$ws.Range("N77").Value = -62721
# Row 80
$ws.Range("H80").Value = 26920
$ws.Range("J80").Value = 26920
$ws.Range("L80").Value = 26920
$ws.Range("N80").Value = -29166
# Row 83
$ws.Range("H83").Value = 26920
$ws.Range("J83").Value = 26920
$ws.Range("L83").Value = 80760
$ws.Range("N83").Value = -91992
# Row 134
$ws.Range("H134").Value = 779.4054
$ws.Range("I134").Value = 720.93335
$ws.Range("J134").Value = 1030
$ws.Range("K134").Value = 2162.80005
$ws.Range("L134").Value = 3090
$ws.Range("M134").Value = 372.1999500000002
$ws.Range("N134").Value = -8160
# Row 138
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 723.4761999999999
$ws.Range("I98").Value = 598.6667
$ws.Range("J98").Value = 817.0833
$ws.Range("K98").Value = 1796.0001
$ws.Range("L98").Value = 2451.2499
$ws.Range("M98").Value = -298.0001
$ws.Range("N98").Value = -5447.2499
# Row 107
$ws.Range("H107").Value = 540.129
$ws.Range("I107").Value = 323.22223
$ws.Range("J107").Value = 840.46155
$ws.Range("K107").Value = 969.66669
$ws.Range("L107").Value = 2521.38465
$ws.Range("M107").Value = 950.33331
$ws.Range("N107").Value = -6361.38465
# Row 113
$ws.Range("H113").Value = 933.62103
$ws.Range("I113").Value = 512.5
$ws.Range("J113").Value = 994.50604
$ws.Range("K113").Value = 1537.5
$ws.Range("L113").Value = 2983.51812
$ws.Range("M113").Value = 632.5
$ws.Range("N113").Value = -7323.51812
# Row 122
$ws.Range("H122").Value = 567.0571
$ws.Range("I122").Value = 411.82758
$ws.Range("J122").Value = 1317.3334
$ws.Range("K122").Value = 3706.44822
$ws.Range("L122").Value = 11856.0006
$ws.Range("M122").Value = -1256.44822
$ws.Range("N122").Value = -16756.0006
# Row 131
$ws.Range("H131").Value = 12501089
$ws.Range("I131").Value = 1938.3334
$ws.Range("J131").Value = 14706822
$ws.Range("K131").Value = 5815.0002
$ws.Range("L131").Value = 44120466
$ws.Range("M131").Value = -775.0002000000004
$ws.Range("N131").Value = -44130546

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 19985.354
$ws.Range("J93").Value = 20859.438
$ws.Range("L93").Value = 20859.438
$ws.Range("N93").Value = -24603.438
# Row 122
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5050
$ws.Range("N122").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 136
$ws.Range("H136").Value = 3433.6155
$ws.Range("I136").Value = 2095.84
$ws.Range("J136").Value = 5822.5
$ws.Range("K136").Value = 6287.52
$ws.Range("L136").Value = 17467.5
$ws.Range("M136").Value = -3737.52
$ws.Range("N136").Value = -22567.5

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 28305536
$ws.Range("I132").Value = 44119170
$ws.Range("J132").Value = 7458.316
$ws.Range("K132").Value = 132357510
$ws.Range("L132").Value = 22374.948
$ws.Range("M132").Value = -132354980
$ws.Range("N132").Value = -27434.948
# Row 136
$ws.Range("H136").Value = 1457.6538
$ws.Range("I136").Value = 994.35297
$ws.Range("J136").Value = 2332.7778
$ws.Range("K136").Value = 2983.05891
$ws.Range("L136").Value = 6998.3334
$ws.Range("M136").Value = -433.0589100000002
$ws.Range("N136").Value = -12098.3334
